# Planner biocodegaming.xlsx - "fixed login and changed everthing to PHP"
#
# 1) Add a daily log entry (column D) + note (column H) for row 17
#    (date 2017-05-15, serial 42884) on the "Projects Overview" sheet.
# 2) Mark that date cell (C17) with the "Good" (green) cell style, same
#    as the other days that already have entries.
# 3) Give column A an explicit width.
# 4) Update the saved view state (selection + scroll position).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Projects Overview")

# --- New log / note text for row 17 ------------------------------------
$ws.Range("D17").Value = "alex heeft een basic forum en damian is met de login en reg bezig"
$ws.Range("H17").Value = "forum wordt ge implementeerd "

# --- Highlight C17 with the built-in "Good" look (green fill / dark green text) ---
# (Set font/fill directly instead of ".Style = ""Good""" so the existing
#  date NumberFormat on the cell is left untouched.)
$c17 = $ws.Range("C17")
$c17.Font.Size = 11
$c17.Font.Color = 24832        # RGB(0,97,0)   -> FF006100
$c17.Interior.Color = 13561798 # RGB(198,239,206) -> FFC6EFCE

# --- Column A width ------------------------------------------------------
$ws.Columns("A").ColumnWidth = 10.285714285714286   # renders as width="11"

# --- View / selection state ----------------------------------------------
$ws.Activate()
$ws.Range("C11").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
